$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description text (column B) for the three product rows
$ws.Range("B2").Value = "one of the best laptop you can get right now with high end specs"
$ws.Range("B3").Value = "Best watch in market with all your daily tracking"
$ws.Range("B4").Value = "best anc earbuds with this price point "

# Adjust row heights: row 2 -> 30, rows 3 & 4 -> default (auto) row height
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

# Move the active selection to B6
$ws.Range("B6").Select()
